# Apply the Jun 4, 2021 16:00 snapshot columns (Daily/Weekly/Monthly/Closing for
# several "as on" timestamps) to the right of the existing data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give the 20 new header cells (row 1) the same bold/bordered/centered style as the
# rest of the header row before filling in their text.
$ws.Range("A1").Copy()
$ws.Range("DC1:DV1").PasteSpecial(-4122)

# New header labels, columns DC1:DV1
$ws.Range("DC1").Value = "Daily as on Jun 4, 16:00"
$ws.Range("DD1").Value = "Weekly as on Jun 4, 16:00"
$ws.Range("DE1").Value = "Monthly as on Jun 4, 16:00"
$ws.Range("DF1").Value = "Closing as on Jun 4, 16:00"
$ws.Range("DG1").Value = "Daily as on Jun 4, 15:59"
$ws.Range("DH1").Value = "Weekly as on Jun 4, 15:59"
$ws.Range("DI1").Value = "Monthly as on Jun 4, 15:59"
$ws.Range("DJ1").Value = "Closing as on Jun 4, 15:59"
$ws.Range("DK1").Value = "Daily as on Jun 4, 15:57"
$ws.Range("DL1").Value = "Weekly as on Jun 4, 15:57"
$ws.Range("DM1").Value = "Monthly as on Jun 4, 15:57"
$ws.Range("DN1").Value = "Closing as on Jun 4, 15:57"
$ws.Range("DO1").Value = "Daily as on Jun 4, 15:58"
$ws.Range("DP1").Value = "Weekly as on Jun 4, 15:58"
$ws.Range("DQ1").Value = "Monthly as on Jun 4, 15:58"
$ws.Range("DR1").Value = "Closing as on Jun 4, 15:58"
$ws.Range("DS1").Value = "Daily as on Jun 4, 15:56"
$ws.Range("DT1").Value = "Weekly as on Jun 4, 15:56"
$ws.Range("DU1").Value = "Monthly as on Jun 4, 15:56"
$ws.Range("DV1").Value = "Closing as on Jun 4, 15:56"

# New data values per stock row (only the timestamp columns that actually have new
# data for that stock are populated, matching the source update).
# Row 2
$ws.Range("DC2").Value = 72.71
$ws.Range("DD2").Value = 61.45
$ws.Range("DE2").Value = 65.53
$ws.Range("DF2").Value = 2190.5
# Row 3
$ws.Range("DC3").Value = 64.93
$ws.Range("DD3").Value = 70.2
$ws.Range("DE3").Value = 66.98
$ws.Range("DF3").Value = 334.95
# Row 4
$ws.Range("DC4").Value = 69.73
$ws.Range("DD4").Value = 65.52
$ws.Range("DE4").Value = 69.49
$ws.Range("DF4").Value = 5993.45
# Row 5
$ws.Range("DC5").Value = 72.66
$ws.Range("DD5").Value = 74.47
$ws.Range("DE5").Value = 71.24
$ws.Range("DF5").Value = 12155.65
# Row 6
$ws.Range("DG6").Value = 56.02
$ws.Range("DH6").Value = 59.05
$ws.Range("DI6").Value = 63.19
$ws.Range("DJ6").Value = 1500.95
# Row 7
$ws.Range("DC7").Value = 56.23
$ws.Range("DD7").Value = 53.86
$ws.Range("DE7").Value = 62.16
$ws.Range("DF7").Value = 680.9
# Row 8
$ws.Range("DG8").Value = 73.62
$ws.Range("DH8").Value = 59.72
$ws.Range("DI8").Value = 61.62
$ws.Range("DJ8").Value = 3067.6
# Row 9
$ws.Range("DC9").Value = 69.4
$ws.Range("DD9").Value = 54.69
$ws.Range("DE9").Value = 53.82
$ws.Range("DF9").Value = 7214.7
# Row 10
$ws.Range("DC10").Value = 68.45
$ws.Range("DD10").Value = 70.1
$ws.Range("DE10").Value = 68.2
$ws.Range("DF10").Value = 4250.05
# Row 11
$ws.Range("DC11").Value = 58.42
$ws.Range("DD11").Value = 63.02
$ws.Range("DE11").Value = 63.47
$ws.Range("DF11").Value = 59.8
# Row 12
$ws.Range("DC12").Value = 69.07
$ws.Range("DD12").Value = 68.03
$ws.Range("DE12").Value = 56.36
$ws.Range("DF12").Value = 377
# Row 13
$ws.Range("DO13").Value = 56.83
$ws.Range("DP13").Value = 66.22
$ws.Range("DQ13").Value = 69.09
$ws.Range("DR13").Value = 224.4
# Row 14
$ws.Range("DC14").Value = 54.51
$ws.Range("DD14").Value = 69.48
$ws.Range("DE14").Value = 64.09
$ws.Range("DF14").Value = 2248.95
# Row 15
$ws.Range("DC15").Value = 47.71
$ws.Range("DD15").Value = 49.87
$ws.Range("DE15").Value = 58.52
$ws.Range("DF15").Value = 532.7
# Row 16
$ws.Range("DC16").Value = 65.28
$ws.Range("DD16").Value = 61.5
$ws.Range("DE16").Value = 45.95
$ws.Range("DF16").Value = 153.1
# Row 17
$ws.Range("DO17").Value = 69.8
$ws.Range("DP17").Value = 76.37
$ws.Range("DQ17").Value = 83.71
$ws.Range("DR17").Value = 541.2
# Row 18
$ws.Range("DC18").Value = 83.43
$ws.Range("DD18").Value = 64.92
$ws.Range("DE18").Value = 45.01
$ws.Range("DF18").Value = 261.35
# Row 19
$ws.Range("DC19").Value = 70.05
$ws.Range("DD19").Value = 51.3
$ws.Range("DE19").Value = 42.41
$ws.Range("DF19").Value = 9.75
# Row 20
$ws.Range("DK20").Value = 58.37
$ws.Range("DL20").Value = 58
$ws.Range("DM20").Value = 50.24
$ws.Range("DN20").Value = 1009.3
# Row 21
$ws.Range("DC21").Value = 75.03
$ws.Range("DD21").Value = 69.19
$ws.Range("DE21").Value = 53.06
$ws.Range("DF21").Value = 114.6
# Row 22
$ws.Range("DC22").Value = 64.74
$ws.Range("DD22").Value = 63.54
$ws.Range("DE22").Value = 69.51
$ws.Range("DF22").Value = 1921.7
# Row 23
$ws.Range("DC23").Value = 48.49
$ws.Range("DD23").Value = 50.81
$ws.Range("DE23").Value = 47.11
$ws.Range("DF23").Value = 208.75
# Row 24
$ws.Range("DC24").Value = 73.46
$ws.Range("DD24").Value = 67.24
$ws.Range("DE24").Value = 62.94
$ws.Range("DF24").Value = 1537.55
# Row 25
$ws.Range("DC25").Value = 56.24
$ws.Range("DD25").Value = 48.11
$ws.Range("DE25").Value = 48.38
$ws.Range("DF25").Value = 168.15
# Row 26
$ws.Range("DC26").Value = 58.36
$ws.Range("DD26").Value = 52.41
$ws.Range("DE26").Value = 54.25
$ws.Range("DF26").Value = 162.45
# Row 27
$ws.Range("DC27").Value = 81.72
$ws.Range("DD27").Value = 69.68
$ws.Range("DE27").Value = 69.54
$ws.Range("DF27").Value = 1509.35
# Row 28
$ws.Range("DC28").Value = 63.78
$ws.Range("DD28").Value = 57.96
$ws.Range("DE28").Value = 56.78
$ws.Range("DF28").Value = 122.2
# Row 29
$ws.Range("DC29").Value = 60.26
$ws.Range("DD29").Value = 65.12
$ws.Range("DE29").Value = 61.95
$ws.Range("DF29").Value = 989.15
# Row 30
$ws.Range("DC30").Value = 75
$ws.Range("DD30").Value = 74.22
$ws.Range("DE30").Value = 67.3
$ws.Range("DF30").Value = 433.6
# Row 31
$ws.Range("DC31").Value = 58.33
$ws.Range("DD31").Value = 67.96
$ws.Range("DE31").Value = 71.03
$ws.Range("DF31").Value = 108.15
# Row 32
$ws.Range("DC32").Value = 53.44
$ws.Range("DD32").Value = 57.15
$ws.Range("DE32").Value = 68.46
$ws.Range("DF32").Value = 3143.75
# Row 33
$ws.Range("DS33").Value = 60.29
$ws.Range("DT33").Value = 61.83
$ws.Range("DU33").Value = 61.78
$ws.Range("DV33").Value = 746.15
# Row 34
$ws.Range("DC34").Value = 63.11
$ws.Range("DD34").Value = 65.26
$ws.Range("DE34").Value = 51.63
$ws.Range("DF34").Value = 10.45
# Row 35
$ws.Range("DC35").Value = 66.15
$ws.Range("DD35").Value = 69.63
$ws.Range("DE35").Value = 72.66
$ws.Range("DF35").Value = 2924.9
# Row 36
$ws.Range("DC36").Value = 68.51
$ws.Range("DD36").Value = 71.02
$ws.Range("DE36").Value = 81.92
$ws.Range("DF36").Value = 833.85
# Row 37
$ws.Range("DC37").Value = 49.86
$ws.Range("DD37").Value = 63.64
$ws.Range("DE37").Value = 66.84
$ws.Range("DF37").Value = 673.95
# Row 38
$ws.Range("DC38").Value = 57.23
$ws.Range("DD38").Value = 77.76
$ws.Range("DE38").Value = 79.15
$ws.Range("DF38").Value = 1120.7
# Row 39
$ws.Range("DC39").Value = 46.72
$ws.Range("DD39").Value = 67.82
$ws.Range("DE39").Value = 83.04
$ws.Range("DF39").Value = 1731.65
